# Sprint 5 - finaliza transacciones: se recorta el set de datos de prueba
# de "generar_clave" a un unico caso (fila 2) en la hoja "Datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# El caso que se conserva (fila 2) deja de traer valores de
# codigoTransaccion/codigoError/resultadoEsperado/numeroDocumento/usuario/
# clave/segundaClave/tipoDocumento: solo quedan idCaso (A2) y orientacion (B2).
$ws.Range("C2:J2").ClearContents()

# Se eliminan los casos 2 y 3 (filas 3 y 4), que ya no se usan.
$ws.Rows("3:4").Delete()

# Deja la celda seleccionada donde quedo el cursor al terminar la edicion.
$ws.Activate() | Out-Null
$ws.Range("L14").Select() | Out-Null
